$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename the sheet
$ws.Name = "Memory map"

# Recolor the merged "mlp_head_dense2_out[0][id]" cell (B41:B44) from the
# green accent fill to the light grey fill used by the rest of that row,
# preserving each row's own border (top/middle/bottom of the merged box).
$rngB = $ws.Range("B41:B44")
$rngB.Interior.ThemeColor = 4
$rngB.Interior.TintAndShade = 0

# Add the new "Memory" label in BP49
$ws.Range("BP49").Value = "Memory"

# Update the current selection / scroll position
$ws.Range("N37").Select()
